$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates in column A (rows 3-21) use DD-MM-YYYY textual format. Excel's COM
# layer would otherwise try to parse strings such as "01-08-2022" as a date
# (US month-day order) and store a date serial number instead of the literal
# text. Force the range to a text number format first so the new values are
# preserved exactly as typed, matching the original inline-string cell type.
$ws.Range("A3:A21").NumberFormat = "@"

# Row 3: 28/07/2022 -> 28-07-2022 ; D3 0->1 ; G3 0->1
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4: 01/08/2022 -> 01-08-2022 ; D4 0->1 ; E4 0->1 ; H4 1->0
$ws.Range("A4").Value = "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5: 04/08/2022 -> 04-08-2022 ; D5 0->1 ; E5 0->1 ; H5 1->0
$ws.Range("A5").Value = "04-08-2022"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 6: 08/08/2022 -> 08-08-2022 (date format only)
$ws.Range("A6").Value = "08-08-2022"

# Row 7: 11/08/2022 -> 11-08-2022 (date format only)
$ws.Range("A7").Value = "11-08-2022"

# Row 8: 15/08/2022 -> 15-08-2022 (date format only)
$ws.Range("A8").Value = "15-08-2022"

# Row 9: 18/08/2022 -> 18-08-2022 (date format only)
$ws.Range("A9").Value = "18-08-2022"

# Row 10: 22/08/2022 -> 22-08-2022 (date format only)
$ws.Range("A10").Value = "22-08-2022"

# Row 11: 25/08/2022 -> 25-08-2022 (date format only)
$ws.Range("A11").Value = "25-08-2022"

# Row 12: 29/08/2022 -> 29-08-2022 ; D12 0->1 ; E12 0->1 ; H12 1->0
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

# Row 13: 01/09/2022 -> 01-09-2022 ; D13 0->1 ; E13 0->1 ; H13 1->0
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

# Row 14: 05/09/2022 -> 05-09-2022 (date format only)
$ws.Range("A14").Value = "05-09-2022"

# Row 15: 08/09/2022 -> 08-09-2022 (date format only)
$ws.Range("A15").Value = "08-09-2022"

# Row 16: 12/09/2022 -> 12-09-2022 (date format only)
$ws.Range("A16").Value = "12-09-2022"

# Row 17: 15/09/2022 -> 15-09-2022 (date format only)
$ws.Range("A17").Value = "15-09-2022"

# Row 18: 19/09/2022 -> 19-09-2022 (date format only)
$ws.Range("A18").Value = "19-09-2022"

# Row 19: 22/09/2022 -> 22-09-2022 (date format only)
$ws.Range("A19").Value = "22-09-2022"

# Row 20: 26/09/2022 -> 26-09-2022 (date format only)
$ws.Range("A20").Value = "26-09-2022"

# Row 21: 29/09/2022 -> 29-09-2022 (date format only)
$ws.Range("A21").Value = "29-09-2022"
